$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates per commit "[UPD] Added new data"

# max_rain_10min_mm: 0.2 -> 0.6
$ws.Range("J2").Value = 0.6

# max_rain_rate_mmh: 0 -> 11.2
$ws.Range("M2").Value = 11.2

# year_min (for mean_pressure_hPa max): 2025 -> 2023 (keep as text, like other year columns)
$ws.Range("AI2").NumberFormat = "@"
$ws.Range("AI2").Value = "2023"

# mean_pressure_hPa: 1020.062541335979 -> 1021.385952264581
$ws.Range("AK2").Value = 1021.385952264581

# year_min (for pcp_acum_month_mm): 2025 -> 2014 (keep as text, like other year columns)
$ws.Range("AL2").NumberFormat = "@"
$ws.Range("AL2").Value = "2014"

# pcp_acum_month_mm: 0.6000000000000001 -> 6
$ws.Range("AN2").Value = 6
